$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2222
$ws1.Range("F4").Value = 95
$ws1.Range("F5").Value = 13328
$ws1.Range("F8").Value = 520
$ws1.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202403/72UDlNoh1711680247000.jpeg"
$ws1.Range("F9").Value = 486
$ws1.Range("F11").Value = 1000
$ws1.Range("F12").Value = 13818
$ws1.Range("F13").Value = 14473
$ws1.Range("F22").Value = 1106
$ws1.Range("F23").Value = 115
$ws1.Range("F25").Value = 5517
$ws1.Range("F27").Value = 973
$ws1.Range("F30").Value = 110

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2222
$ws4.Range("F4").Value = 95
$ws4.Range("F5").Value = 13329
$ws4.Range("F9").Value = 520
$ws4.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202403/72UDlNoh1711680247000.jpeg"
$ws4.Range("F10").Value = 486
$ws4.Range("F12").Value = 1000
$ws4.Range("F13").Value = 13818
$ws4.Range("F14").Value = 14473
$ws4.Range("F23").Value = 1106
$ws4.Range("F24").Value = 115
$ws4.Range("F26").Value = 5517
$ws4.Range("F28").Value = 973
$ws4.Range("F31").Value = 110
